$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores numeric-looking values as TEXT (e.g. "1.000",
# "27.846.33" -- note the double "thousands" dots, not valid numbers). Force
# the cells we are about to rewrite to Text format first so Excel keeps the
# new values as text instead of auto-converting them to numbers.
$ws.Range("D2:D15").NumberFormat = "@"
$ws.Range("D17:D19").NumberFormat = "@"
$ws.Range("D23:D24").NumberFormat = "@"
$ws.Range("D26:D35").NumberFormat = "@"
$ws.Range("D37:D39").NumberFormat = "@"
$ws.Range("D41:D51").NumberFormat = "@"

# Update Price (column D) cells
$ws.Range("D2").Value = "27.859.76"
$ws.Range("D3").Value = "1.906.88"
$ws.Range("D4").Value = "0.9997"
$ws.Range("D5").Value = "313.30"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").Value = "0.5007"
$ws.Range("D8").Value = "0.3812"
$ws.Range("D9").Value = "0.07278"
$ws.Range("D10").Value = "0.9105"
$ws.Range("D11").Value = "20.95"
$ws.Range("D12").Value = "0.07648"
$ws.Range("D13").Value = "1.885.11"
$ws.Range("D14").Value = "5.498"
$ws.Range("D15").Value = "91.96"
$ws.Range("D17").Value = "0.000008742"
$ws.Range("D18").Value = "0.9994"
$ws.Range("D19").Value = "27.889.67"
$ws.Range("D23").Value = "6.579"
$ws.Range("D24").Value = "153.20"
$ws.Range("D26").Value = "2.218"
$ws.Range("D27").Value = "18.39"
$ws.Range("D28").Value = "115.44"
$ws.Range("D29").Value = "4.921"
$ws.Range("D30").Value = "0.09030"
$ws.Range("D31").Value = "3.193"
$ws.Range("D32").Value = "1.229"
$ws.Range("D33").Value = "4.779"
$ws.Range("D34").Value = "0.7735"
$ws.Range("D35").Value = "0.02086"
$ws.Range("D37").Value = "1.094"
$ws.Range("D38").Value = "0.5572"
$ws.Range("D39").Value = "3.021"
$ws.Range("D41").Value = "6.904"
$ws.Range("D42").Value = "8.486"
$ws.Range("D43").Value = "0.1516"
$ws.Range("D44").Value = "111.77"
$ws.Range("D45").Value = "0.4837"
$ws.Range("D46").Value = "10.60"
$ws.Range("D47").Value = "1.000"
$ws.Range("D48").Value = "1.636"
$ws.Range("D49").Value = "67.54"
$ws.Range("D50").Value = "0.06057"
$ws.Range("D51").Value = "0.9090"

# Update Volume(1h) (column E) cells
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("E22").Value = "  -0.61%  "
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("E26").Value = "  +3.90%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("E36").Value = "  -3.40%  "
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("E49").Value = "  -0.52%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  +1.10%  "
